$d = $word.ActiveDocument
$sec = $d.Sections(1)
$f = $sec.Footers(2)
$r = $f.Range
$ok = $r.Find.Execute("52", $true, $false, $false, $false, $false, $true, 1, $false, "76", 2)
Write-Output ("find result=" + $ok)
Write-Output ("after text=[" + $f.Range.Text + "]")
